$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.621.57'
$ws.Range("E2").Value = '  +0.51%  '

$ws.Range("D3").Value = '3.683.82'
$ws.Range("E3").Value = '  +0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '669.05'
$ws.Range("E5").Value = '  -0.97%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.89'
$ws.Range("E6").Value = '  +1.28%  '

$ws.Range("E7").Value = '  +0.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.501'
$ws.Range("E8").Value = '  +1.54%  '

$ws.Range("E9").Value = '  +0.11%  '

$ws.Range("E10").Value = '  +2.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.441'
$ws.Range("E11").Value = '  +1.35%  '

$ws.Range("E12").Value = '  +0.90%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '33.02'
$ws.Range("E13").Value = '  +2.22%  '

$ws.Range("D14").Value = '3.691.37'
$ws.Range("E14").Value = '  +0.78%  '

$ws.Range("D15").Value = '69.609.47'
$ws.Range("E15").Value = '  +0.56%  '

$ws.Range("E16").Value = '  +2.48%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '16.13'
$ws.Range("E17").Value = '  +0.45%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.46'
$ws.Range("E18").Value = '  +0.48%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '470.88'
$ws.Range("E19").Value = '  +0.78%  '

$ws.Range("E20").Value = '  -2.49%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.646'
$ws.Range("E21").Value = '  -0.38%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '79.76'
$ws.Range("E22").Value = '  +0.05%  '

$ws.Range("D23").Value = '3.832.87'
$ws.Range("E23").Value = '  +0.37%  '

$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.03%  '

$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000127'
$ws.Range("E25").Value = '  +4.00%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.94'
$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.05'
$ws.Range("E27").Value = '  +0.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.68'
$ws.Range("E28").Value = '  +0.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.70'
$ws.Range("E29").Value = '  -2.11%  '

$ws.Range("E30").Value = '  +1.27%  '

$ws.Range("E31").Value = '  +0.20%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.166'
$ws.Range("E32").Value = '  +3.32%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.75'
$ws.Range("E33").Value = '  -0.35%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.49'
$ws.Range("E34").Value = '  -1.96%  '

$ws.Range("D35").Value = '3.689.09'
$ws.Range("E35").Value = '  +0.59%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.46'
$ws.Range("E36").Value = '  +3.80%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.09'
$ws.Range("E37").Value = '  -2.98%  '

$ws.Range("E38").Value = '  -0.06%  '

$ws.Range("E39").Value = '  +1.60%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  -0.06%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '177.00'
$ws.Range("E41").Value = '  +1.31%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0906'
$ws.Range("E42").Value = '  +0.83%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.934'
$ws.Range("E43").Value = '  -0.58%  '

$ws.Range("E44").Value = '  -1.14%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.74'
$ws.Range("E45").Value = '  +1.55%  '

$ws.Range("E46").Value = '  +0.47%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000271'
$ws.Range("E47").Value = '  -1.88%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.28'
$ws.Range("E48").Value = '  -2.59%  '

$ws.Range("E49").Value = '  -0.25%  '

$ws.Range("E50").Value = '  +1.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '364.45'
$ws.Range("E51").Value = '  +1.04%  '
